$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-12: test case data (fills in numbers/text first) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Player throws rocks, tries to run, is caught by the bear, and ends up at the bear's cabin"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Player throws rocks, tries to run, is caught by the bear, but escapes the forest by taking the right path"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "2 (Tree height: 20)"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "Player throws rocks, climbs a tall tree, avoids the bear, but ends up at the bear's cabin"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "2 (Tree height: 20)"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Player throws rocks, climbs a tall tree, avoids the bear, and escapes the forest by taking the right path"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "2 (Tree height: 10)"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Player throws rocks, climbs a short tree, is caught by the bear, and the game ends"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Speed: 25"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "Player runs quickly, outruns the bear, but ends up at the bear's cabin"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Speed: 25"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Player runs quickly, outruns the bear, and escapes the forest by taking the right path"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "Speed: 15"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = "Player tries to run but is caught by the bear, and the game ends"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "N/A"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "Player stays still, bear loses interest, but ends up at the bear's cabin after taking the left path"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "N/A"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "Player stays still, bear loses interest, and escapes the forest by taking the right path"

$ws.Range("A3:E12").HorizontalAlignment = -4152   # xlRight

# --- Row 2: header row, right aligned + bold ---
$ws.Range("A2").Value = "Test Case"
$ws.Range("B2").Value = "Choice1"
$ws.Range("C2").Value = "Choice2 (Ft/MPH)"
$ws.Range("D2").Value = "Choice3"
$ws.Range("E2").Value = "Expected Output"
$ws.Range("A2:E2").HorizontalAlignment = -4152   # xlRight
$ws.Range("A2:E2").Font.Bold = $true

# --- Row 1: shrink merge, drop trailing style on E1, right-align A1:D1 ---
$ws.Range("A1:E1").UnMerge()
$ws.Range("E1").ClearFormats()
$ws.Range("A1:D1").HorizontalAlignment = -4152   # xlRight
$ws.Range("A1:D1").Merge()

# --- Column widths (approximate autofit results baked into the source file) ---
$ws.Columns("A").ColumnWidth = 8.166666666666666
$ws.Columns("B").ColumnWidth = 9.830729166666666
$ws.Columns("C").ColumnWidth = 18.330729166666668
$ws.Columns("D").ColumnWidth = 13.666666666666666
$ws.Columns("E").ColumnWidth = 83.33072916666667

# --- Selection, matches the saved view state ---
$ws.Range("E19").Select()
